$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old group labels to new group labels.
$map = @{
    "Group 6"  = "Group 3"
    "Group 2"  = "Group 6"
    "Group 3"  = "Group 7"
    "Group 1"  = "Group 5"
    "Group 5"  = "Group 8"
    "Group 11" = "Group 11"
    "Group 8"  = "Group 9"
}

$used = $ws.UsedRange
$startRow = $used.Row
$startCol = $used.Column
$endRow = $startRow + $used.Rows.Count - 1
$endCol = $startCol + $used.Columns.Count - 1

for ($r = $startRow; $r -le $endRow; $r++) {
    for ($c = $startCol; $c -le $endCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $map.ContainsKey([string]$val)) {
            $cell.Value2 = $map[[string]$val]
        }
    }
}
